$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 18288.45
$ws.Range("I15").Value = 18288.45
$ws.Range("K15").Value = 54865.35000000001
$ws.Range("M15").Value = -54696.35000000001
$ws.Range("H18").Value = 1278.6
$ws.Range("I18").Value = 1278.6
$ws.Range("K18").Value = 1278.6
$ws.Range("M18").Value = -994.5999999999999
$ws.Range("H113").Value = 16668500
$ws.Range("I113").Value = 50001000
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 50001000
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = -49997746
$ws.Range("N113").Value = -8758
$ws.Range("H132").Value = 2710985.5
$ws.Range("I132").Value = 999.8205
$ws.Range("J132").Value = 55555704
$ws.Range("K132").Value = 2999.4615
$ws.Range("L132").Value = 166667112
$ws.Range("M132").Value = -469.4615000000003
$ws.Range("N132").Value = -166672172
$ws.Range("H141").Value = 2850.4707
$ws.Range("I141").Value = 1271.4166
$ws.Range("J141").Value = 6640.2
$ws.Range("K141").Value = 3814.2498
$ws.Range("L141").Value = 19920.6
$ws.Range("M141").Value = 1365.7502
$ws.Range("N141").Value = -30280.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17838.682
$ws.Range("I2").Value = 19259.21
$ws.Range("K2").Value = 19259.21
$ws.Range("M2").Value = -19146.21
$ws.Range("H45").Value = 417889.72
$ws.Range("I45").Value = 715267.0600000001
$ws.Range("K45").Value = 715267.0600000001
$ws.Range("M45").Value = -714890.0600000001
$ws.Range("H63").Value = 2933.3333
$ws.Range("I63").Value = 2033.3334
$ws.Range("J63").Value = 4733.3335
$ws.Range("K63").Value = 2033.3334
$ws.Range("L63").Value = 4733.3335
$ws.Range("M63").Value = -1347.3334
$ws.Range("N63").Value = -6105.3335
$ws.Range("H66").Value = 2933.3333
$ws.Range("I66").Value = 2033.3334
$ws.Range("J66").Value = 4733.3335
$ws.Range("K66").Value = 10166.667
$ws.Range("L66").Value = 23666.6675
$ws.Range("M66").Value = -6734.666999999999
$ws.Range("N66").Value = -30530.6675
$ws.Range("H74").Value = 44148308
$ws.Range("I74").Value = 42238760
$ws.Range("J74").Value = 48488196
$ws.Range("K74").Value = 42238760
$ws.Range("L74").Value = 48488196
$ws.Range("M74").Value = -42237886
$ws.Range("N74").Value = -48489944
$ws.Range("H77").Value = 44148308
$ws.Range("I77").Value = 42238760
$ws.Range("J77").Value = 48488196
$ws.Range("K77").Value = 211193800
$ws.Range("L77").Value = 242440980
$ws.Range("M77").Value = -211189432
$ws.Range("N77").Value = -242449716
$ws.Range("H116").Value = 17838.682
$ws.Range("I116").Value = 19259.21
$ws.Range("K116").Value = 19259.21
$ws.Range("M116").Value = -16965.21
$ws.Range("H122").Value = 1827.9584
$ws.Range("I122").Value = 1349.5883
$ws.Range("K122").Value = 4048.7649
$ws.Range("M122").Value = -1598.7649
$ws.Range("H123").Value = 49916.332
$ws.Range("J123").Value = 49916.332
$ws.Range("L123").Value = 49916.332
$ws.Range("N123").Value = -59716.332
$ws.Range("H132").Value = 8010654.5
$ws.Range("I132").Value = 5557540.5
$ws.Range("J132").Value = 21298356
$ws.Range("K132").Value = 16672621.5
$ws.Range("L132").Value = 63895068
$ws.Range("M132").Value = -16670091.5
$ws.Range("N132").Value = -63900128

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17838.682
$ws.Range("I3").Value = 19259.21
$ws.Range("K3").Value = 19259.21
$ws.Range("M3").Value = -19145.21
$ws.Range("H35").Value = 34980
$ws.Range("J35").Value = 34980
$ws.Range("L35").Value = 34980
$ws.Range("N35").Value = -35600
$ws.Range("H82").Value = 23345.562
$ws.Range("I82").Value = 2337.7144
$ws.Range("J82").Value = 39685
$ws.Range("K82").Value = 2337.7144
$ws.Range("L82").Value = 39685
$ws.Range("M82").Value = -1954.7144
$ws.Range("N82").Value = -40451
$ws.Range("H85").Value = 23345.562
$ws.Range("I85").Value = 2337.7144
$ws.Range("J85").Value = 39685
$ws.Range("K85").Value = 2337.7144
$ws.Range("L85").Value = 39685
$ws.Range("M85").Value = -1011.7144
$ws.Range("N85").Value = -42337
$ws.Range("H105").Value = 1584.3438
$ws.Range("I105").Value = 1581.6666
$ws.Range("K105").Value = 1581.6666
$ws.Range("M105").Value = 165.3334
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120
$ws.Range("H134").Value = 12312663
$ws.Range("I134").Value = 12195872
$ws.Range("K134").Value = 36587616
$ws.Range("M134").Value = -36585081

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10527886
$ws.Range("I31").Value = 25642112
$ws.Range("J31").Value = 1906.8214
$ws.Range("K31").Value = 25642112
$ws.Range("L31").Value = 1906.8214
$ws.Range("M31").Value = -25641817
$ws.Range("N31").Value = -2496.8214
$ws.Range("H34").Value = 10527886
$ws.Range("I34").Value = 25642112
$ws.Range("J34").Value = 1906.8214
$ws.Range("K34").Value = 25642112
$ws.Range("L34").Value = 1906.8214
$ws.Range("M34").Value = -25641910
$ws.Range("N34").Value = -2310.8214
$ws.Range("H94").Value = 18522668
$ws.Range("I94").Value = 1327.75
$ws.Range("J94").Value = 21743772
$ws.Range("K94").Value = 1327.75
$ws.Range("L94").Value = 21743772
$ws.Range("M94").Value = -876.75
$ws.Range("N94").Value = -21744674

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2321.2266
$ws.Range("I113").Value = 2731.6667
$ws.Range("J113").Value = 1705.5667
$ws.Range("K113").Value = 8195.000100000001
$ws.Range("L113").Value = 5116.7001
$ws.Range("M113").Value = -6025.000100000001
$ws.Range("N113").Value = -9456.7001
$ws.Range("H123").Value = 62501700
$ws.Range("I123").Value = 83333940
$ws.Range("J123").Value = 4999
$ws.Range("K123").Value = 250001820
$ws.Range("L123").Value = 14997
$ws.Range("M123").Value = -249999370
$ws.Range("N123").Value = -19897
$ws.Range("H125").Value = 47620904
$ws.Range("I125").Value = 333334530
$ws.Range("J125").Value = 1967.1666
$ws.Range("K125").Value = 1000003590
$ws.Range("L125").Value = 5901.4998
$ws.Range("M125").Value = -999998670
$ws.Range("N125").Value = -15741.4998
$ws.Range("H129").Value = 5320848
$ws.Range("I129").Value = 19232042
$ws.Range("J129").Value = 1862.3823
$ws.Range("K129").Value = 57696126
$ws.Range("L129").Value = 5587.1469
$ws.Range("M129").Value = -57691126
$ws.Range("N129").Value = -15587.1469
$ws.Range("H130").Value = 1283.3334
$ws.Range("I130").Value = 1200
$ws.Range("J130").Value = 1450
$ws.Range("K130").Value = 3600
$ws.Range("L130").Value = 4350
$ws.Range("M130").Value = 1420
$ws.Range("N130").Value = -14390
$ws.Range("H131").Value = 42598.293
$ws.Range("I131").Value = 111585.555
$ws.Range("J131").Value = 1205.9333
$ws.Range("K131").Value = 334756.665
$ws.Range("L131").Value = 3617.7999
$ws.Range("M131").Value = -329716.665
$ws.Range("N131").Value = -13697.7999
$ws.Range("H139").Value = 57137.332
$ws.Range("J139").Value = 5033
$ws.Range("L139").Value = 15099
$ws.Range("N139").Value = -25379

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H132").Value = 8282727.5
$ws.Range("I132").Value = 5349479.5
$ws.Range("J132").Value = 17985010
$ws.Range("K132").Value = 16048438.5
$ws.Range("L132").Value = 53955030
$ws.Range("M132").Value = -16045908.5
$ws.Range("N132").Value = -53960090

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 23744.666
$ws.Range("I107").Value = 29700.285
$ws.Range("J107").Value = 2900
$ws.Range("K107").Value = 89100.855
$ws.Range("L107").Value = 8700
$ws.Range("M107").Value = -87180.855
$ws.Range("N107").Value = -12540
